# Revert "Merge pull request #18 ... US1023501_TINStatusPages"
# Restores the original (pre-merge) CSR test-case row and re-adds the
# "Create Enrollment" test-case row on Sheet2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Restore the fully-qualified class name for the CSR View Payments test.
$ws.Range("C10").Value = "test.java.TestCSRViewPayments"

# Re-add the "Create Enrollment" test case as row 11.
$ws.Range("A11").Value = "UPA_Regression"
$ws.Range("B11").Value = "Create Enrollment"
$ws.Range("C11").Value = "test.java.TestCreateEnrollment"

# Keep the active selection in sync with the newly added last row.
$ws.Range("C11").Select()
